$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$ws.Range("A1").Value = [double]"0.789663034731252"
$ws.Range("B1").Value = [double]"0.052572570399218"
$ws.Range("C1").Value = [double]"2.83056368339569e-11"
$ws.Range("D1").Value = [double]"0.0526011057121666"
$ws.Range("E1").Value = [double]"2.65968735340929e-11"
$ws.Range("F1").Value = [double]"3.00160935843614e-11"
$ws.Range("G1").Value = [double]"0.0525806972260075"
$ws.Range("H1").Value = [double]"0.0525825918464371"
$ws.Range("A2").Value = [double]"0.0525735342504167"
$ws.Range("B2").Value = [double]"0.789724479023232"
$ws.Range("C2").Value = [double]"0.052576006684229"
$ws.Range("D2").Value = [double]"1.55755044729494e-10"
$ws.Range("E2").Value = [double]"0.0525648588550806"
$ws.Range("F2").Value = [double]"0.05256112097702"
$ws.Range("G2").Value = [double]"2.73248396790522e-11"
$ws.Range("H2").Value = [double]"2.6941730299439e-11"
$ws.Range("A3").Value = [double]"0.285687175291914"
$ws.Range("B3").Value = [double]"1.42779935995239e-10"
$ws.Range("C3").Value = [double]"0.714312823850704"
$ws.Range("D3").Value = [double]"1.43056191175965e-10"
$ws.Range("E3").Value = [double]"1.43100277461616e-10"
$ws.Range("F3").Value = [double]"1.4288612786737e-10"
$ws.Range("G3").Value = [double]"1.42779606799385e-10"
$ws.Range("H3").Value = [double]"1.42779996888038e-10"
$ws.Range("A4").Value = [double]"1.42736474335368e-10"
$ws.Range("B4").Value = [double]"0.285664487067984"
$ws.Range("C4").Value = [double]"1.42736474338538e-10"
$ws.Range("D4").Value = [double]"0.714335512075074"
$ws.Range("E4").Value = [double]"1.42736474370231e-10"
$ws.Range("F4").Value = [double]"1.42736474329029e-10"
$ws.Range("G4").Value = [double]"1.42953414099726e-10"
$ws.Range("H4").Value = [double]"1.43043205301378e-10"
$ws.Range("A5").Value = [double]"0.117647584506745"
$ws.Range("B5").Value = [double]"5.8816207924013e-11"
$ws.Range("C5").Value = [double]"2.80865499152992e-10"
$ws.Range("D5").Value = [double]"5.88389336381614e-11"
$ws.Range("E5").Value = [double]"0.882352414560867"
$ws.Range("F5").Value = [double]"4.16233700271662e-10"
$ws.Range("G5").Value = [double]"5.88147002427473e-11"
$ws.Range("H5").Value = [double]"5.88183941092504e-11"
$ws.Range("A6").Value = [double]"0.117656457796993"
$ws.Range("B6").Value = [double]"5.88290301908951e-11"
$ws.Range("C6").Value = [double]"5.92350192718063e-11"
$ws.Range("D6").Value = [double]"5.88185323070693e-11"
$ws.Range("E6").Value = [double]"5.92545882550221e-11"
$ws.Range("F6").Value = [double]"0.88234354184922"
$ws.Range("G6").Value = [double]"5.88210852726765e-11"
$ws.Range("H6").Value = [double]"5.88294711562214e-11"
$ws.Range("A7").Value = [double]"5.87994279675011e-11"
$ws.Range("B7").Value = [double]"0.117632073477369"
$ws.Range("C7").Value = [double]"5.87994279675011e-11"
$ws.Range("D7").Value = [double]"5.88061114325301e-11"
$ws.Range("E7").Value = [double]"5.87994279675011e-11"
$ws.Range("F7").Value = [double]"5.87994279675011e-11"
$ws.Range("G7").Value = [double]"0.882367926165822"
$ws.Range("H7").Value = [double]"6.28058114034742e-11"
$ws.Range("A8").Value = [double]"5.88245202978735e-11"
$ws.Range("B8").Value = [double]"0.117653367610261"
$ws.Range("C8").Value = [double]"5.88138340234195e-11"
$ws.Range("D8").Value = [double]"5.89792940070867e-06"
$ws.Range("E8").Value = [double]"5.88138340064425e-11"
$ws.Range("F8").Value = [double]"5.88450114436511e-11"
$ws.Range("G8").Value = [double]"6.55885766211661e-11"
$ws.Range("H8").Value = [double]"0.882340734159453"
$ws.Range("D13").Select() | Out-Null
